$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: date header (mirrors row 1: merged A:C, style s4) ---
# Stage the literal text "06/21/2023" in a scratch cell via a formula so it
# lands as a plain string (not an auto-parsed date serial), then bring just
# the value across before applying row 1's formatting + merge.
$ws.Range("Z1").Formula = "=T(""06/21/2023"")"
$ws.Range("Z1").Copy()
$ws.Range("A11").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("A1:C1").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122)
$ws.Range("A11:C11").Merge()

# --- Row 12: column headers (mirrors row 2: style s2) ---
$ws.Range("A2:C2").Copy()
$ws.Range("A12:C12").PasteSpecial(-4122)
$ws.Range("A12").Value = "TASKS"
$ws.Range("B12").Value = "ASSIGNED TO"
$ws.Range("C12").Value = "PROGRESS"

# --- Rows 13-19: data rows (mirror row 3: style s3) ---
$ws.Range("A3:C3").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Range("A16:C16").PasteSpecial(-4122)
$ws.Range("A17:C17").PasteSpecial(-4122)
$ws.Range("A18:C18").PasteSpecial(-4122)
$ws.Range("A19:C19").PasteSpecial(-4122)

$ws.Range("A13").Value = "Database"
$ws.Range("B13").Value = "Platon, Forbes, Rodriguez"
$ws.Range("C13").Value = "Completed"

$ws.Range("A14").Value = "Design Presentation"
$ws.Range("B14").Value = "Platon, Forbes, Rodriguez"
$ws.Range("C14").Value = "Completed"

$ws.Range("A15").Value = "Homepage"
$ws.Range("B15").Value = "Platon"
$ws.Range("C15").Value = "Completed"

$ws.Range("A16").Value = "BP Forms"
$ws.Range("B16").Value = "Forbes"
$ws.Range("C16").Value = "Completed"

$ws.Range("A17").Value = "Admin Homepage"
$ws.Range("B17").Value = "Rodriguez"
$ws.Range("C17").Value = "In-progress"

$ws.Range("A18").Value = "Admin Review Boards"
$ws.Range("B18").Value = "Forbes"
$ws.Range("C18").Value = "In-progress"

$ws.Range("A19").Value = "Documentation (push to repo)"
$ws.Range("B19").Value = "Rodriguez"
$ws.Range("C19").Value = "In-progress"

# --- Selection matches the target view state ---
$ws.Range("G12").Select()
